$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows whose C / D / E values change (names stay line1..line6 for rows 2-7) ---

# Row 8 (index 6): name becomes "line7" (was "extr1"), C 5->14, D 12->11, E 0->1(true)
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9 (index 7): name becomes "line8" (was "extr2"), C 5->16, D stays 9, E stays true
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Row 10 (index 8): name becomes "extr1" (was "extr3"), C 10->5, D 11->12, E stays true
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 11 (index 9): name becomes "extr2" (was "extr4"), C 7->5, D 8->9, E stays true
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 12 (index 10): name becomes "extr3" (was "extr5"), C 9->10, D stays 11, E true->false
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

# Row 13 (index 11): name becomes "extr4" (was "extr6"), C stays 7, D 11->8, E true->false
$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

# Row 14 (index 12): name becomes "extr5" (was "extr7"), C 5->9, D 7->11, E true->false
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# Row 15 (index 13): name becomes "extr6" (was "extr8"), C 8->7, D 5->11, E false->true
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# --- New rows 16 and 17 ---

# Row 16 (index 14): name "extr7", C 5, D 7, E true
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

# Row 17 (index 15): name "extr8", C 8, D 5, E true
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true

# Column A on rows 2-15 uses a bold/bordered/centered style (s="1"); replicate it
# on the two new rows by copying the format from an existing styled cell.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
